$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the now-unused trailing rows (13-16) -----------------------
$ws.Rows("13:16").Delete()

# --- New columns G (CRITICAL CODE) and H (CRITICAL NAME) ---------------
$ws.Columns("G").ColumnWidth = 15.166666666666666
$ws.Columns("H").ColumnWidth = 23.166666666666668

# --- Header row ----------------------------------------------------------
$ws.Range("G1").NumberFormat = "@"
$ws.Range("G1").Value = "CRITICAL CODE"
$ws.Range("H1").NumberFormat = "@"
$ws.Range("H1").Value = "CRITICAL NAME"

# --- Row 2: single student/topic instead of comma-joined pair ----------
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "19110373"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "Pham Quang Hung"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "TL-10"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "HUDHJDHJ"

# --- Reshuffled TOPIC CODE (column C) values for rows 3-12 --------------
$topicCodes = @{
    3  = "TL-9"
    4  = "TL-4"
    5  = "TL-8"
    6  = "TL-5"
    7  = "TL-2"
    8  = "TL-3"
    9  = "TL-12"
    10 = "TL-6"
    11 = "TL-14"
    12 = "TL-11"
}
foreach ($r in $topicCodes.Keys) {
    $cell = $ws.Range("C$r")
    $cell.NumberFormat = "@"
    $cell.Value = $topicCodes[$r]
}

# --- New CRITICAL CODE / CRITICAL NAME columns for every data row ------
for ($r = 2; $r -le 12; $r++) {
    $gCell = $ws.Range("G$r")
    $gCell.NumberFormat = "@"
    $gCell.Value = "2000"

    $hCell = $ws.Range("H$r")
    $hCell.NumberFormat = "@"
    $hCell.Value = "LECTURER HUNG"
}
